# Applies the "Updated symbol list" crypto price/volume refresh to sheet1.
# Values are assigned with a leading apostrophe to force text storage,
# matching the original inlineStr (text) cell type used for Price/Volume columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'279.59"
$ws.Range("E2").Value = "'0.87%"

# Row 3
$ws.Range("D3").Value = "'27.48"
$ws.Range("E3").Value = "'1.20%"

# Row 4
$ws.Range("D4").Value = "'4.834"
$ws.Range("E4").Value = "'-2.28%"

# Row 5
$ws.Range("D5").Value = "'0.06403"
$ws.Range("E5").Value = "'-0.15%"

# Row 6
$ws.Range("D6").Value = "'7.039"
$ws.Range("E6").Value = "'1.68%"

# Row 7
$ws.Range("D7").Value = "'1.328"
$ws.Range("E7").Value = "'6.17%"

# Row 8
$ws.Range("D8").Value = "'0.9021"
$ws.Range("E8").Value = "'2.28%"

# Row 9
$ws.Range("D9").Value = "'0.1538"
$ws.Range("E9").Value = "'0.83%"

# Row 10
$ws.Range("D10").Value = "'0.06134"
$ws.Range("E10").Value = "'20.86%"

# Row 11
$ws.Range("D11").Value = "'0.07463"
$ws.Range("E11").Value = "'-1.13%"

# Row 12
$ws.Range("D12").Value = "'0.02923"
$ws.Range("E12").Value = "'1.24%"

# Row 13
$ws.Range("D13").Value = "'0.08992"
$ws.Range("E13").Value = "'-0.23%"

# Row 14
$ws.Range("D14").Value = "'0.001597"
$ws.Range("E14").Value = "'2.02%"

# Row 15
$ws.Range("D15").Value = "'0.0006435"
$ws.Range("E15").Value = "'0.18%"

# Row 16
$ws.Range("D16").Value = "'0.006025"
$ws.Range("E16").Value = "'0.62%"

# Row 17
$ws.Range("D17").Value = "'3.489"
$ws.Range("E17").Value = "'0.99%"

# Row 18
$ws.Range("D18").Value = "'3.306"
$ws.Range("E18").Value = "'-0.48%"

# Row 19
$ws.Range("D19").Value = "'2.231"
$ws.Range("E19").Value = "'-1.81%"

# Row 21
$ws.Range("D21").Value = "'0.1353"
$ws.Range("E21").Value = "'1.21%"

# Row 22
$ws.Range("D22").Value = "'3.906"
$ws.Range("E22").Value = "'-0.07%"

# Row 23
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1504"
$ws.Range("E23").Value = "'8.94%"

# Row 24
$ws.Range("B24").Value = "'CoinExToken"
$ws.Range("C24").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04396"
$ws.Range("E24").Value = "'-0.81%"

# Row 25
$ws.Range("D25").Value = "'0.001174"
$ws.Range("E25").Value = "'-0.09%"

# Row 26
$ws.Range("D26").Value = "'0.004292"
$ws.Range("E26").Value = "'10.92%"

# Row 28
$ws.Range("D28").Value = "'0.0001178"
$ws.Range("E28").Value = "'-1.96%"

# Row 29
$ws.Range("D29").Value = "'0.0001656"
$ws.Range("E29").Value = "'-14.48%"

# Row 40
$ws.Range("D40").Value = "'0.04078"
$ws.Range("E40").Value = "'-1.53%"

# Row 41
$ws.Range("D41").Value = "'0.006589"
$ws.Range("E41").Value = "'-3.00%"

# Row 42
$ws.Range("D42").Value = "'0.1399"
$ws.Range("E42").Value = "'18.60%"

# Row 43
$ws.Range("D43").Value = "'0.002087"
$ws.Range("E43").Value = "'-9.01%"

# Row 44
$ws.Range("D44").Value = "'0.01100"
$ws.Range("E44").Value = "'-2.27%"

# Row 45
$ws.Range("D45").Value = "'0.00005547"
$ws.Range("E45").Value = "'6.50%"

# Row 46
$ws.Range("E46").Value = "'9.82%"

# Row 47
$ws.Range("D47").Value = "'0.01847"
$ws.Range("E47").Value = "'-8.89%"
